$d = $word.ActiveDocument

# =========================================================================
# Change 1: "Login eller oprettelse af bruger skal udføres." becomes three
# runs reading "Login eller o" | "prettelse af bruger skal behandles
# fortroligt" | "." (net text: "...skal behandles fortroligt.")
# =========================================================================
$i = 0
$paraIdx = -1
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Login eller oprettelse*") {
        $paraIdx = $i
    }
}

$p = $d.Paragraphs($paraIdx)
$pStart = $p.Range.Start

$prefix = "Login eller o"
$boundary1 = $pStart + $prefix.Length

# Insert the new middle-run text right after "Login eller o"; toggling a
# character attribute forces Word to keep it as its own run once the
# attribute is reverted.
$newMiddle = "prettelse af bruger skal behandles fortroligt"
$insPoint = $d.Range($boundary1, $boundary1)
$insPoint.InsertAfter($newMiddle)
$middleRange = $d.Range($boundary1, $boundary1 + $newMiddle.Length)
$middleRange.Bold = 1
$middleRange.Bold = 0

# Remove the stale tail text ("prettelse af bruger skal udføres"), but keep
# its trailing period in place (it becomes its own run) so its original
# formatting (da-DK language) is preserved.
$oldTail = "prettelse af bruger skal udføres."
$oldTailStart = $boundary1 + $newMiddle.Length
$oldTailDeleteLen = $oldTail.Length - 1
$deleteRange = $d.Range($oldTailStart, $oldTailStart + $oldTailDeleteLen)
$deleteRange.Text = ""

$dotRange = $d.Range($oldTailStart, $oldTailStart + 1)
$dotRange.Bold = 1
$dotRange.Bold = 0

# =========================================================================
# Change 2: insert two new list paragraphs ("Historik" at level 2 and
# "Oversigt over historik med eksportmulighed af en CSV-fil anvendes." at
# level 3) right before the "Usability" paragraph.
# =========================================================================
$i = 0
$usabilityIdx = -1
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "Usability*") {
        $usabilityIdx = $i
    }
}

$target = $d.Paragraphs($usabilityIdx)
$target.Range.InsertParagraphBefore()
$historikRange = $d.Paragraphs($usabilityIdx).Range
$historikRange.Text = "Historik"
$historikRange.ListFormat.ListLevelNumber = 2

$usabilityIdx2 = $usabilityIdx + 1
$target2 = $d.Paragraphs($usabilityIdx2)
$target2.Range.InsertParagraphBefore()
$oversigtRange = $d.Paragraphs($usabilityIdx2).Range
$oversigtRange.Text = "Oversigt over historik med eksportmulighed af en CSV-fil anvendes."
$oversigtRange.ListFormat.ListLevelNumber = 3

# =========================================================================
# Change 3: remove the stray "_GoBack" bookmark left around "Supportability".
# =========================================================================
if ($d.Bookmarks.Exists("_GoBack")) {
    $bm = $d.Bookmarks("_GoBack")
    $bm.Delete()
}
